$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update 최종점수 (K column) - each row decreases by 0.5
$ws.Range("K2").Value = 58.7
$ws.Range("K3").Value = 54.7
$ws.Range("K4").Value = 52.7
$ws.Range("K5").Value = 51.7

# Update MACRO_SCORE (N column) - all rows updated to the same new value
$ws.Range("N2:N5").Value = 49.16024380385575
